# Update "Update latest output (run 9)" changes
$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 1534.7215845
$schedule.Range("F2").Value = 25.37568757440476

# --- Sheet: Detailed ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B11").Value = 65.15832
$detailed.Range("B12").Value = 81.07834

$detailed.Range("B13").Value = 79.35364
$detailed.Range("C13").Value = "historical"

$detailed.Range("B14").Value = 80.02
$detailed.Range("C14").Value = "historical"

$detailed.Range("B15").Value = 63.99078
$detailed.Range("B16").Value = 36.07
$detailed.Range("B17").Value = 36.07
$detailed.Range("B18").Value = 36.06
$detailed.Range("B19").Value = 26.42341

$detailed.Range("B21").Value = 0.51
$detailed.Range("B22").Value = 0.51
$detailed.Range("B23").Value = 5.21834

$detailed.Range("B25").Value = 0.51

$detailed.Range("B28").Value = 0.51

$detailed.Range("B30").Value = 53.42659

$detailed.Range("B32").Value = 58.80779
$detailed.Range("B33").Value = 56.98
$detailed.Range("B34").Value = 56.69206
$detailed.Range("B35").Value = 61.2163
$detailed.Range("B36").Value = 66.20182
$detailed.Range("B37").Value = 32.55525
$detailed.Range("B38").Value = 52.09869
$detailed.Range("B39").Value = 73.69302

$detailed.Range("B41").Value = 120.01
$detailed.Range("B42").Value = 158.99
$detailed.Range("B43").Value = 120.01
$detailed.Range("B44").Value = 105.79
$detailed.Range("B45").Value = 105.79
$detailed.Range("B46").Value = 120.01
$detailed.Range("B47").Value = 69.11084
$detailed.Range("B48").Value = 58.17358
$detailed.Range("B49").Value = 60.55376
